# This edit reorders the data rows (rows 2-15) of the worksheet: the same
# 14 records are redistributed across the rows (no records are added or
# removed), matching the "Automatic update of files" re-fetch from source.
#
# Mapping of new row -> old row (source of its data):
#   2<-5  3<-2  4<-6  5<-7  6<-8  7<-9  8<-3  9<-10  10<-11
#   11<-12  12<-13  13<-14  14<-15  15<-4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that actually vary between the data records and therefore need to
# be moved along with the row. (Other columns such as C,D,S,T,U,V,W,Y,Z,AA,
# AB,AD,AE,AG,AW,AX are identical across every data row so they do not need
# to be touched.)
$cols = @("A", "B", "E", "F", "G", "H", "P", "Q", "R")

$firstDataRow = 2
$lastDataRow = 15

# Snapshot the current values of every relevant cell before we start
# overwriting anything.
$snapshot = @{}
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# New row -> old row mapping (where the data for the new row comes from).
$rowMap = @{
    2  = 5
    3  = 2
    4  = 6
    5  = 7
    6  = 8
    7  = 9
    8  = 3
    9  = 10
    10 = 11
    11 = 12
    12 = 13
    13 = 14
    14 = 15
    15 = 4
}

foreach ($newRow in ($rowMap.Keys | Sort-Object)) {
    $oldRow = $rowMap[$newRow]
    $src = $snapshot[$oldRow]
    foreach ($c in $cols) {
        $ws.Range("$c$newRow").Value = $src[$c]
    }
}
